$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new trade row (row 9) to the CELG random trade log.
$ws.Range("A9").Value = 42654.743784722225
$ws.Range("A9").NumberFormat = "m/d/yy h:mm"

$ws.Range("B9").Value = $true
$ws.Range("C9").Value = 10249.81
$ws.Range("D9").Value = 10242.129999999999
$ws.Range("E9").Value = 104.43
$ws.Range("F9").Value = 104.269997

$ws.Range("G9").Value = $true
$ws.Range("G9").NumberFormat = "m/d/yy h:mm"

$ws.Range("H9").Value = -0.15
$ws.Range("I9").Value = $false

# Widen column A slightly to fit the new content (bestFit recalculation).
# (The runtime's char-width -> stored-width rounding means 15.375 itself is not
# reproducible; 14.5 is the input that lands closest, at 15.33(3).)
$ws.Columns("A").ColumnWidth = 14.5
